$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Capture the "latest row" date style (date-only, no time) from row 24
# before it gets reassigned to the regular datetime style.
$latestDateFormat = $ws.Range("A24").NumberFormat

# Row 24 was previously the "latest" row and used a special date-only
# number format. Since row 25 becomes the new latest row, row 24 reverts
# to the regular datetime format used by every other non-final row.
$ws.Range("A24").NumberFormat = $ws.Range("A2").NumberFormat

# Append the new day's data as the new "latest" row.
$ws.Range("A25").Value = 45609
$ws.Range("B25").Value = 63
$ws.Range("C25").Value = 53
$ws.Range("D25").Value = 59

# Give the new latest row's date cell the special date-only format.
$ws.Range("A25").NumberFormat = $latestDateFormat
